# Automatische test-sync: 2025-07-31 22:06:50
#
# Adds a new log entry (row 22) to the "Logs" sheet, adds the matching
# aggregated category row (row 8) to the "Dashboard" sheet, and extends
# the bar chart series + conditional formatting ranges so they include
# the new rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 22
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A22").Value = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("B22").Value = "mailmind.test@zohomail.eu"
$logs.Range("C22").Value = "Testmail #20: Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("D22").Value = "Klacht / Probleem"
$logs.Range("E22").Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$logs.Range("F22").Value = "2025-07-31 22:06:18"
$logs.Range("G22").Value = "Ja"
$logs.Range("H22").Value = "Ja"
$logs.Range("I22").Value = "Nee"
$logs.Range("J22").Value = "Nee"

# Extend the conditional-formatting ranges on the Logs sheet from row 21
# to row 22 for every formatted column.
$logsCfColumns = "D", "G", "H", "I", "J"
foreach ($col in $logsCfColumns) {
    $oldRange = $logs.Range($col + "2:" + $col + "21")
    $newRange = $logs.Range($col + "2:" + $col + "22")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append row 8 with the new category total
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A8").Value = "Klacht / Probleem"
$dash.Range("B8").Value = 1

# ---------------------------------------------------------------------
# 3. Update the bar chart so its category/value series cover A2:A8/B2:B8
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$8,'Dashboard'!`$B`$2:`$B`$8,1)"
